$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "57.965.51"
$ws.Range("E2").Value = "  +1.92%  "
# Row 3
$ws.Range("D3").Value = "2.344.18"
$ws.Range("E3").Value = "  +1.17%  "
# Row 4
$ws.Range("E4").Value = "  -0.14%  "
# Row 5
$ws.Range("D5").Value = "539.62"
$ws.Range("E5").Value = "  +1.93%  "
# Row 6
$ws.Range("D6").Value = "135.32"
$ws.Range("E6").Value = "  +2.44%  "
# Row 7
$ws.Range("E7").Value = "  +0.51%  "
# Row 8
$ws.Range("D8").Value = "0.566"
$ws.Range("E8").Value = "  +5.95%  "
# Row 9
$ws.Range("E9").Value = "  +0.37%  "
# Row 10
$ws.Range("E10").Value = "  +4.80%  "
# Row 11
$ws.Range("E11").Value = "  -0.67%  "
# Row 12
$ws.Range("E12").Value = "  +1.71%  "
# Row 13
$ws.Range("D13").Value = "23.75"
$ws.Range("E13").Value = "  +1.24%  "
# Row 14
$ws.Range("D14").Value = "2.761.80"
$ws.Range("E14").Value = "  +0.92%  "
# Row 15
$ws.Range("D15").Value = "57.937.65"
$ws.Range("E15").Value = "  +1.78%  "
# Row 16
$ws.Range("E16").Value = "  +0.54%  "
# Row 17
$ws.Range("D17").Value = "2.358.03"
$ws.Range("E17").Value = "  +0.76%  "
# Row 18
$ws.Range("D18").Value = "10.68"
$ws.Range("E18").Value = "  +2.49%  "
# Row 19
$ws.Range("D19").Value = "331.36"
$ws.Range("E19").Value = "  -1.50%  "
# Row 20
$ws.Range("E20").Value = "  +2.62%  "
# Row 21
$ws.Range("D21").Value = "6.76"
$ws.Range("E21").Value = "  -1.33%  "
# Row 22
$ws.Range("E22").Value = "  -0.01%  "
# Row 23
$ws.Range("D23").Value = "62.65"
# Row 24
$ws.Range("E24").Value = "  -0.35%  "
# Row 25
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").Value = "8.48"
$ws.Range("E25").Value = "  -2.59%  "
# Row 26
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  +0.50%  "
# Row 27
$ws.Range("E27").Value = "  +1.76%  "
# Row 28
$ws.Range("E28").Value = "  +1.61%  "
# Row 29
$ws.Range("D29").Value = "171.42"
$ws.Range("E29").Value = "  -0.73%  "
# Row 30
$ws.Range("D30").Value = "0.0₃0735"
$ws.Range("E30").Value = "  +1.56%  "
# Row 31
$ws.Range("D31").Value = "6.12"
$ws.Range("E31").Value = "  +0.27%  "
# Row 32
$ws.Range("E32").Value = "  +12.11%  "
# Row 33
$ws.Range("D33").Value = "18.42"
# Row 35
$ws.Range("B35").Value = "FirstDigitalUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.61%  "
# Row 36
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Value = "4.22"
$ws.Range("E36").Value = "  +6.05%  "
# Row 37
$ws.Range("E37").Value = "  -0.52%  "
# Row 39
$ws.Range("E39").Value = "  -0.05%  "
# Row 40
$ws.Range("D40").Value = "145.07"
$ws.Range("E40").Value = "  -2.63%  "
# Row 41
$ws.Range("D41").Value = "292.49"
$ws.Range("E41").Value = "  +3.51%  "
# Row 42
$ws.Range("E42").Value = "  +0.59%  "
# Row 43
$ws.Range("E43").Value = "  +0.68%  "
# Row 44
$ws.Range("D44").Value = "0.0947"
$ws.Range("E44").Value = "  +1.82%  "
# Row 45
$ws.Range("D45").Value = "19.18"
$ws.Range("E45").Value = "  +1.81%  "
# Row 46
$ws.Range("E46").Value = "  +0.39%  "
# Row 47
$ws.Range("E47").Value = "  +0.58%  "
# Row 48
$ws.Range("E48").Value = "  +1.28%  "
# Row 49
$ws.Range("D49").Value = "0.381"
$ws.Range("E49").Value = "  -0.14%  "
# Row 50
$ws.Range("D50").Value = "17.46"
$ws.Range("E50").Value = "  +0.08%  "
# Row 51
$ws.Range("D51").Value = "11.08"
$ws.Range("E51").Value = "  +0.47%  "
